$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text first so numeric-looking strings
# (e.g. "212.76") are not silently coerced into floating point numbers,
# which would corrupt values like "1.90" or "0.0591".
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.559.96'
$ws.Range('E2').Value = '  +0.86%  '

$ws.Range('D3').Value = '1.569.85'
$ws.Range('E3').Value = '  -1.41%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Value = '212.76'
$ws.Range('E5').Value = '  -0.59%  '

$ws.Range('E6').Value = '  -0.42%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').Value = '45.79'
$ws.Range('E8').Value = '  +4.16%  '

$ws.Range('D9').Value = '24.09'
$ws.Range('E9').Value = '  +0.15%  '

$ws.Range('E10').Value = '  -1.74%  '

$ws.Range('D11').Value = '0.0591'
$ws.Range('E11').Value = '  -1.62%  '

$ws.Range('D12').Value = '0.0888'
$ws.Range('E12').Value = '  -0.14%  '

$ws.Range('D13').Value = '1.792.93'
$ws.Range('E13').Value = '  -1.54%  '

$ws.Range('D14').Value = '1.555.21'
$ws.Range('E14').Value = '  -1.73%  '

$ws.Range('E15').Value = '  -1.98%  '

$ws.Range('D16').Value = '28.534.05'
$ws.Range('E16').Value = '  +0.70%  '

$ws.Range('E17').Value = '  -2.07%  '

$ws.Range('D18').Value = '62.23'
$ws.Range('E18').Value = '  -1.39%  '

$ws.Range('D19').Value = '230.17'
$ws.Range('E19').Value = '  +1.19%  '

$ws.Range('D20').Value = '7.36'
$ws.Range('E20').Value = '  -1.47%  '

$ws.Range('E21').Value = '  -2.62%  '

$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('D23').Value = '3.87'
$ws.Range('E23').Value = '  -5.91%  '

$ws.Range('D24').Value = '9.12'
$ws.Range('E24').Value = '  -2.22%  '

$ws.Range('D25').Value = '2.13'
$ws.Range('E25').Value = '  +8.92%  '

$ws.Range('D26').Value = '151.52'
$ws.Range('E26').Value = '  -0.04%  '

$ws.Range('D27').Value = '15.03'
$ws.Range('E27').Value = '  -1.11%  '

$ws.Range('E28').Value = '  -2.48%  '

$ws.Range('E29').Value = '  -3.23%  '

$ws.Range('E30').Value = '  -0.05%  '

$ws.Range('D31').Value = '0.0484'
$ws.Range('E31').Value = '  +2.12%  '

$ws.Range('D32').Value = '1.11'
$ws.Range('E32').Value = '  -3.05%  '

$ws.Range('E33').Value = '  -1.01%  '

$ws.Range('E34').Value = '  -1.11%  '

$ws.Range('D35').Value = '1.392.99'
$ws.Range('E35').Value = '  -0.43%  '

$ws.Range('E36').Value = '  +0.77%  '

$ws.Range('D37').Value = '1.53'
$ws.Range('E37').Value = '  -3.75%  '

$ws.Range('D38').Value = '2.37'
$ws.Range('E38').Value = '  +0.68%  '

$ws.Range('E39').Value = '  +2.44%  '

$ws.Range('E40').Value = '  -0.70%  '

$ws.Range('E41').Value = '  -3.22%  '

$ws.Range('E42').Value = '  +0.02%  '

$ws.Range('D43').Value = '1.90'
$ws.Range('E43').Value = '  +0.33%  '

$ws.Range('E44').Value = '  -3.01%  '

$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').Value = '0.0463'
$ws.Range('E45').Value = '  +0.47%  '

$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '5.49'
$ws.Range('E46').Value = '  -2.96%  '

$ws.Range('D47').Value = '0.971'
$ws.Range('E47').Value = '  -1.41%  '

$ws.Range('D48').Value = '62.86'
$ws.Range('E48').Value = '  -2.18%  '

$ws.Range('D49').Value = '1.705.52'
$ws.Range('E49').Value = '  -1.60%  '

$ws.Range('D50').Value = '86.29'
$ws.Range('E50').Value = '  -1.41%  '

$ws.Range('E51').Value = '  -0.23%  '

# Restore the original (default) cell style for the Price column now
# that the text values are safely stored, so formatting matches the
# source workbook exactly.
$ws.Range("D2:D51").Style = "Normal"
